$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "Name"
$ws.Range("C2").Select()
